# Update gh-pages output data for 苏州-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2-18
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 758
$ws1.Range("F3").Value = 672
$ws1.Range("F4").Value = 1189
$ws1.Range("G5").Value = 50
$ws1.Range("F6").Value = 605
$ws1.Range("G6").Value = 70
$ws1.Range("G7").Value = 75
$ws1.Range("F13").Value = 301
$ws1.Range("F17").Value = 11142
$ws1.Range("F18").Value = 5325

# Sheet "全部类型" (All types) - rows 2-21
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 758
$ws4.Range("F3").Value = 672
$ws4.Range("F4").Value = 1189
$ws4.Range("G5").Value = 50
$ws4.Range("F6").Value = 605
$ws4.Range("G6").Value = 70
$ws4.Range("G7").Value = 75
$ws4.Range("F15").Value = 301
$ws4.Range("F19").Value = 11142
$ws4.Range("F21").Value = 5325
